$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the "word-guidance-*" snippet ids back to "word-scenarios-*"
# (moving the examples back to their own "scenarios" folder).

# Rows 9-11: "word-guidance-doc-assembly" -> "word-scenarios-doc-assembly"
$ws.Range("C9").Value = "word-scenarios-doc-assembly"
$ws.Range("C10").Value = "word-scenarios-doc-assembly"
$ws.Range("C11").Value = "word-scenarios-doc-assembly"

# Rows 32-33: "word-guidance-multiple-property-set" -> "word-scenarios-multiple-property-set"
$ws.Range("C32").Value = "word-scenarios-multiple-property-set"
$ws.Range("C33").Value = "word-scenarios-multiple-property-set"

# Update the active selection to match the final saved state
$ws.Range("C33").Select()
